# Update "想去人数" (F column) counts on the "展览" (sheet1), "演出" (sheet2),
# and "全部类型" (sheet4) worksheets to reflect newly generated output.

$wb = $excel.ActiveWorkbook

# --- 展览 (Exhibition) sheet ---
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F2").Value = 7036
$wsExhibition.Range("F4").Value = 70
$wsExhibition.Range("F7").Value = 6975
$wsExhibition.Range("F11").Value = 23
$wsExhibition.Range("F17").Value = 52
$wsExhibition.Range("F18").Value = 49
$wsExhibition.Range("F20").Value = 5339
$wsExhibition.Range("F21").Value = 131
$wsExhibition.Range("F22").Value = 186
$wsExhibition.Range("F23").Value = 763
$wsExhibition.Range("F25").Value = 264

# --- 演出 (Performance) sheet ---
$wsPerformance = $wb.Worksheets.Item("演出")
$wsPerformance.Range("F2").Value = 1

# --- 全部类型 (All Types) sheet ---
$wsAllTypes = $wb.Worksheets.Item("全部类型")
$wsAllTypes.Range("F2").Value = 7036
$wsAllTypes.Range("F4").Value = 70
$wsAllTypes.Range("F7").Value = 6975
$wsAllTypes.Range("F11").Value = 23
$wsAllTypes.Range("F17").Value = 52
$wsAllTypes.Range("F18").Value = 49
$wsAllTypes.Range("F20").Value = 1
$wsAllTypes.Range("F21").Value = 5339
$wsAllTypes.Range("F23").Value = 131
$wsAllTypes.Range("F24").Value = 186
$wsAllTypes.Range("F25").Value = 763
$wsAllTypes.Range("F27").Value = 264
